$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-14 19:13:04", 0.0004),
    @("2023-12-14 19:13:46", 0.0026),
    @("2023-12-14 19:14:45", 0.0038),
    @("2023-12-14 19:14:58", 0.0004)
)

$startRow = 335
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
